$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "土地" (Land) — rename headers to machine-readable names
# and append normalized metadata columns I:O
# ============================================================
$wsLand = $wb.Worksheets.Item("土地")

$wsLand.Range("B1").Value = "name"
$wsLand.Range("C1").Value = "area"
$wsLand.Range("D1").Value = "share_portion"
$wsLand.Range("E1").Value = "owner"
$wsLand.Range("F1").Value = "register_date"
$wsLand.Range("G1").Value = "register_reason"
$wsLand.Range("H1").Value = "acquire_value"

# Clean up stray whitespace / OCR artifacts in the existing data
$wsLand.Range("B2").Value = "臺北市大安區仁愛段二小段06010000地號"
$wsLand.Range("D2").Value = "30000分之1703"
$wsLand.Range("F2").Value = "93年01月02日"

$wsLand.Range("B3").Value = "臺北市中山區長安段四小段02980000地號"
$wsLand.Range("D3").Value = "10000分之211"
$wsLand.Range("F3").Value = "92年01月17日"
$wsLand.Range("G3").Value = "1III：■■!■貝買"

$wsLand.Range("B4").Value = "桃園縣平鎮市忠貞段00010003地號"
$wsLand.Range("D4").Value = "10000分之173"
$wsLand.Range("F4").Value = "102年02月20曰"
$wsLand.Range("H4").Value = "14000000(房地總價額）"

# New metadata columns
$wsLand.Range("I1").Value = "property_category"
$wsLand.Range("J1").Value = "category"
$wsLand.Range("K1").Value = "date"
$wsLand.Range("L1").Value = "legislator_name"
$wsLand.Range("M1").Value = "legislator_id"
$wsLand.Range("N1").Value = "source_file"
$wsLand.Range("O1").Value = "index"

$wsLand.Range("I2:I4").Value = "land"
$wsLand.Range("J2:J4").Value = "normal"
$wsLand.Range("K2:K4").NumberFormat = "@"
$wsLand.Range("K2:K4").Value = "2013-07-15"
$wsLand.Range("L2:L4").Value = "林世嘉"
$wsLand.Range("M2:M4").Value = 1740
$wsLand.Range("N2:N4").Value = "tmpf70f1"

$wsLand.Range("O2").Value = 15
$wsLand.Range("O3").Value = 16
$wsLand.Range("O4").Value = 17

# ============================================================
# Sheet "建物" (Building) — same whitespace / OCR cleanup
# ============================================================
$wsBuilding = $wb.Worksheets.Item("建物")

$wsBuilding.Range("B2").Value = "臺北市大安區仁愛段二小段05345000建號"
$wsBuilding.Range("F2").Value = "93年01月02日"

$wsBuilding.Range("B3").Value = "臺北市大安區仁愛段二小段03854000建號"
$wsBuilding.Range("F3").Value = "93年01月02日"
$wsBuilding.Range("G3").Value = "I""rl""T广貝買"

$wsBuilding.Range("B4").Value = "臺北市中山區長安段四小段01678000建號"
$wsBuilding.Range("F4").Value = "92年01月17日"

$wsBuilding.Range("B5").Value = "桃園縣平鎮市忠貞段01899000建號"
$wsBuilding.Range("F5").Value = "102年02月20曰"
$wsBuilding.Range("G5").Value = "rift■■■r:貝買"
$wsBuilding.Range("H5").Value = "14000000(房地總價額）"

# ============================================================
# Sheet "汽車" (Car) — same whitespace cleanup
# ============================================================
$wsCar = $wb.Worksheets.Item("汽車")

$wsCar.Range("B2").Value = "福斯PASSAT2.0"
$wsCar.Range("E2").Value = "94年01月28日"

# ============================================================
# Sheet "債務" (Debt) — same whitespace cleanup
# ============================================================
$wsDebt = $wb.Worksheets.Item("債務")

$wsDebt.Range("D2").Value = "兆豐國際商業銀行台北復興分行"
$wsDebt.Range("F2").Value = "民國97年07月"

$wsDebt.Range("D3").Value = "兆豐國際商業銀行台北復興分行"
$wsDebt.Range("F3").Value = "民國93年02月"

$wsDebt.Range("F4").Value = "民國102年02月"
$wsDebt.Range("F5").Value = "民國98年04月"
$wsDebt.Range("F6").Value = "民國92年01月"
$wsDebt.Range("F7").Value = "民國92年01月"
$wsDebt.Range("F8").Value = "民國92年01月"

$wsDebt.Range("D9").Value = "郭素珍臺北市北投區明德路"
$wsDebt.Range("E9").NumberFormat = "@"
$wsDebt.Range("E9").Value = "2000000"
$wsDebt.Range("F9").Value = "民國92年02月"
